$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41..56 down to 42..57
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly data point
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value = "Ñuble"
$ws.Cells.Item(41, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 100112022
$ws.Cells.Item(41, 7).Value = "Arveja Verde"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 60
$ws.Cells.Item(41, 11).Value = 23000
$ws.Cells.Item(41, 12).Value = 24000
$ws.Cells.Item(41, 13).Value = 23500
$ws.Cells.Item(41, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(41, 16).Value = 940
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
